# fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: advance the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the price list values in column D (rows 29-37)
$ws.Range("D29").Value = 94.935
$ws.Range("D30").Value = 101.646
$ws.Range("D31").Value = 106.442
$ws.Range("D32").Value = 111.716
$ws.Range("D33").Value = 116.993
$ws.Range("D34").Value = 123.699
$ws.Range("D35").Value = 134.252
$ws.Range("D36").Value = 151.514
$ws.Range("D37").Value = 177.405
